$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 4137.273
$ws.Range("I8").Value = 4536
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 13608
$ws.Range("L8").Value = 450
$ws.Range("M8").Value = -13469
$ws.Range("H64").Value = 6255.4443
$ws.Range("I64").Value = 4100
$ws.Range("J64").Value = 6524.875
$ws.Range("K64").Value = 4100
$ws.Range("L64").Value = 6524.875
$ws.Range("M64").Value = -3852
$ws.Range("N64").Value = -7020.875
$ws.Range("H67").Value = 6255.4443
$ws.Range("I67").Value = 4100
$ws.Range("J67").Value = 6524.875
$ws.Range("K67").Value = 4100
$ws.Range("L67").Value = 6524.875
$ws.Range("M67").Value = -3242
$ws.Range("N67").Value = -8240.875
$ws.Range("H137").Value = 52364.7
$ws.Range("I137").Value = 2399.75
$ws.Range("K137").Value = 7199.25
$ws.Range("M137").Value = -4649.25
$ws.Range("H141").Value = 72108.69500000001
$ws.Range("I141").Value = 102984.555
$ws.Range("K141").Value = 308953.665
$ws.Range("M141").Value = -303773.665
$ws.Range("N8").Value = -728

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2153.2173
$ws.Range("I2").Value = 1733.6842
$ws.Range("J2").Value = 4146
$ws.Range("K2").Value = 1733.6842
$ws.Range("L2").Value = 4146
$ws.Range("M2").Value = -1620.6842
$ws.Range("N2").Value = -4372
$ws.Range("H32").Value = 37573.965
$ws.Range("I32").Value = 22302.043
$ws.Range("J32").Value = 117327.336
$ws.Range("K32").Value = 22302.043
$ws.Range("L32").Value = 117327.336
$ws.Range("M32").Value = -22015.043
$ws.Range("N32").Value = -117901.336
$ws.Range("H37").Value = 50000000
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("H103").Value = 29241.334
$ws.Range("J103").Value = 28862
$ws.Range("L103").Value = 28862
$ws.Range("N103").Value = -31206
$ws.Range("H116").Value = 2153.2173
$ws.Range("I116").Value = 1733.6842
$ws.Range("J116").Value = 4146
$ws.Range("K116").Value = 1733.6842
$ws.Range("L116").Value = 4146
$ws.Range("M116").Value = 560.3158000000001
$ws.Range("N116").Value = -8734
$ws.Range("H132").Value = 2265.5908
$ws.Range("I132").Value = 1664.0555
$ws.Range("K132").Value = 4992.166499999999
$ws.Range("M132").Value = -2462.166499999999
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2153.2173
$ws.Range("I3").Value = 1733.6842
$ws.Range("J3").Value = 4146
$ws.Range("K3").Value = 1733.6842
$ws.Range("L3").Value = 4146
$ws.Range("M3").Value = -1619.6842
$ws.Range("N3").Value = -4374
$ws.Range("H134").Value = 2322.5881
$ws.Range("I134").Value = 2190
$ws.Range("J134").Value = 4444
$ws.Range("K134").Value = 6570
$ws.Range("L134").Value = 13332
$ws.Range("M134").Value = -4035
$ws.Range("N134").Value = -18402
$ws.Range("H135").Value = 79996.664
$ws.Range("J135").Value = 79996.664
$ws.Range("L135").Value = 79996.664
$ws.Range("N135").Value = -90136.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4434
$ws.Range("I31").Value = 3044.2856
$ws.Range("J31").Value = 5244.6665
$ws.Range("K31").Value = 3044.2856
$ws.Range("L31").Value = 5244.6665
$ws.Range("M31").Value = -2749.2856
$ws.Range("N31").Value = -5834.6665
$ws.Range("H34").Value = 4434
$ws.Range("I34").Value = 3044.2856
$ws.Range("J34").Value = 5244.6665
$ws.Range("K34").Value = 3044.2856
$ws.Range("L34").Value = 5244.6665
$ws.Range("M34").Value = -2842.2856
$ws.Range("N34").Value = -5648.6665
$ws.Range("H62").Value = 3364.2856
$ws.Range("I62").Value = 3025
$ws.Range("K62").Value = 3025
$ws.Range("M62").Value = -2401
$ws.Range("H65").Value = 3364.2856
$ws.Range("I65").Value = 3025
$ws.Range("K65").Value = 15125
$ws.Range("M65").Value = -12005
$ws.Range("H132").Value = 988.2
$ws.Range("I132").Value = 988.2
$ws.Range("K132").Value = 2964.6
$ws.Range("M132").Value = -434.6000000000004
$ws.Range("H134").Value = 2073.4092
$ws.Range("I134").Value = 1895.6842
$ws.Range("K134").Value = 5687.0526
$ws.Range("M134").Value = -3152.0526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 930.1177
$ws.Range("I5").Value = 554.5
$ws.Range("J5").Value = 1135
$ws.Range("K5").Value = 1663.5
$ws.Range("L5").Value = 3405
$ws.Range("M5").Value = -1551.5
$ws.Range("N5").Value = -3629
$ws.Range("H40").Value = 2349.4614
$ws.Range("J40").Value = 3775.25
$ws.Range("L40").Value = 15101
$ws.Range("N40").Value = -15239
$ws.Range("H135").Value = 930.1177
$ws.Range("I135").Value = 554.5
$ws.Range("J135").Value = 1135
$ws.Range("K135").Value = 4990.5
$ws.Range("L135").Value = 10215
$ws.Range("M135").Value = -2455.5
$ws.Range("N135").Value = -15285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2029.375
$ws.Range("I102").Value = 1407
$ws.Range("J102").Value = 2402.8
$ws.Range("K102").Value = 1407
$ws.Range("L102").Value = 2402.8
$ws.Range("M102").Value = 215
$ws.Range("N102").Value = -5646.8
$ws.Range("H113").Value = 50003300
$ws.Range("I113").Value = 62502624
$ws.Range("K113").Value = 62502624
$ws.Range("M113").Value = -62500454
$ws.Range("H126").Value = 2606.5
$ws.Range("I126").Value = 1823
$ws.Range("K126").Value = 5469
$ws.Range("M126").Value = -2999
$ws.Range("H132").Value = 1779
$ws.Range("I132").Value = 1519.75
$ws.Range("J132").Value = 2297.5
$ws.Range("K132").Value = 4559.25
$ws.Range("L132").Value = 6892.5
$ws.Range("M132").Value = -2029.25
$ws.Range("N132").Value = -11952.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 24112.154
$ws.Range("I7").Value = 32495.445
$ws.Range("J7").Value = 5249.75
$ws.Range("K7").Value = 32495.445
$ws.Range("L7").Value = 5249.75
$ws.Range("M7").Value = -32383.445
$ws.Range("N7").Value = -5473.75
$ws.Range("H40").Value = 6088.1113
$ws.Range("I40").Value = 5266.5
$ws.Range("K40").Value = 5266.5
$ws.Range("M40").Value = -5130.5
$ws.Range("H126").Value = 24112.154
$ws.Range("I126").Value = 32495.445
$ws.Range("J126").Value = 5249.75
$ws.Range("K126").Value = 97486.33499999999
$ws.Range("L126").Value = 15749.25
$ws.Range("M126").Value = -95016.33499999999
$ws.Range("N126").Value = -20689.25
$ws.Range("H132").Value = 2831
$ws.Range("J132").Value = 3267.923
$ws.Range("L132").Value = 9803.769
$ws.Range("N132").Value = -14863.769
$ws.Range("H136").Value = 3285.647
$ws.Range("I136").Value = 2420.375
$ws.Range("J136").Value = 4054.7778
$ws.Range("K136").Value = 7261.125
$ws.Range("L136").Value = 12164.3334
$ws.Range("M136").Value = -4711.125
$ws.Range("N136").Value = -17264.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 40857.332
$ws.Range("J97").Value = 40857.332
$ws.Range("L97").Value = 40857.332
$ws.Range("N97").Value = -42839.332
$ws.Range("H126").Value = 3444.2
$ws.Range("I126").Value = 1740.3334
$ws.Range("K126").Value = 5221.0002
$ws.Range("M126").Value = -2751.0002
$ws.Range("H132").Value = 14380.5
$ws.Range("I132").Value = 17534.2
$ws.Range("K132").Value = 52602.60000000001
$ws.Range("M132").Value = -50072.60000000001
